$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3266.0264
$ws.Range("I15").Value = 3266.0264
$ws.Range("K15").Value = 9798.0792
$ws.Range("M15").Value = -9629.0792
$ws.Range("H28").Value = 671.8333
$ws.Range("I28").Value = 551.25
$ws.Range("J28").Value = 768.3
$ws.Range("K28").Value = 551.25
$ws.Range("L28").Value = 768.3
$ws.Range("M28").Value = -66.25
$ws.Range("N28").Value = -1738.3
$ws.Range("H86").Value = 8174.8423
$ws.Range("I86").Value = 9061.571
$ws.Range("K86").Value = 9061.571
$ws.Range("M86").Value = -7938.571
$ws.Range("H89").Value = 8174.8423
$ws.Range("I89").Value = 9061.571
$ws.Range("K89").Value = 45307.855
$ws.Range("M89").Value = -39691.855
$ws.Range("H92").Value = 1546.6923
$ws.Range("I92").Value = 360.3
$ws.Range("K92").Value = 360.3
$ws.Range("M92").Value = 887.7
$ws.Range("H137").Value = 4168841.8
$ws.Range("I137").Value = 2084942.2
$ws.Range("J137").Value = 8336641
$ws.Range("K137").Value = 6254826.6
$ws.Range("L137").Value = 25009923
$ws.Range("M137").Value = -6252276.6
$ws.Range("N137").Value = -25015023

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2406.889
$ws.Range("I61").Value = 2027
$ws.Range("J61").Value = 3166.6667
$ws.Range("K61").Value = 2027
$ws.Range("L61").Value = 3166.6667
$ws.Range("M61").Value = -1815
$ws.Range("N61").Value = -3590.6667
$ws.Range("H74").Value = 892.6829
$ws.Range("I74").Value = 1027.2593
$ws.Range("J74").Value = 633.1429000000001
$ws.Range("K74").Value = 1027.2593
$ws.Range("L74").Value = 633.1429000000001
$ws.Range("M74").Value = -153.2592999999999
$ws.Range("N74").Value = -2381.1429
$ws.Range("H77").Value = 892.6829
$ws.Range("I77").Value = 1027.2593
$ws.Range("J77").Value = 633.1429000000001
$ws.Range("K77").Value = 5136.296499999999
$ws.Range("L77").Value = 3165.7145
$ws.Range("M77").Value = -768.2964999999995
$ws.Range("N77").Value = -11901.7145
$ws.Range("H132").Value = 210762.5
$ws.Range("I132").Value = 252040.1
$ws.Range("J132").Value = 4374.5
$ws.Range("K132").Value = 756120.3
$ws.Range("L132").Value = 13123.5
$ws.Range("M132").Value = -753590.3
$ws.Range("N132").Value = -18183.5
$ws.Range("H136").Value = 2406.889
$ws.Range("I136").Value = 2027
$ws.Range("J136").Value = 3166.6667
$ws.Range("K136").Value = 6081
$ws.Range("L136").Value = 9500.000100000001
$ws.Range("M136").Value = -3531
$ws.Range("N136").Value = -14600.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1054.65
$ws.Range("I94").Value = 933.5333000000001
$ws.Range("J94").Value = 1418
$ws.Range("K94").Value = 933.5333000000001
$ws.Range("L94").Value = 1418
$ws.Range("M94").Value = -482.5333000000001
$ws.Range("N94").Value = -2320
$ws.Range("H107").Value = 4693.1025
$ws.Range("I107").Value = 4932.4517
$ws.Range("J107").Value = 3765.625
$ws.Range("K107").Value = 4932.4517
$ws.Range("L107").Value = 3765.625
$ws.Range("M107").Value = -3012.4517
$ws.Range("N107").Value = -7605.625
$ws.Range("H134").Value = 303053.6
$ws.Range("I134").Value = 378117
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 1134351
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -1131816
$ws.Range("N134").Value = -13470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 2161.5
$ws.Range("I8").Value = 999
$ws.Range("J8").Value = 2394
$ws.Range("K8").Value = 999
$ws.Range("L8").Value = 2394
$ws.Range("M8").Value = -859
$ws.Range("N8").Value = -2674
$ws.Range("H31").Value = 2306.037
$ws.Range("J31").Value = 2985.8333
$ws.Range("L31").Value = 2985.8333
$ws.Range("N31").Value = -3575.8333
$ws.Range("H34").Value = 2306.037
$ws.Range("J34").Value = 2985.8333
$ws.Range("L34").Value = 2985.8333
$ws.Range("N34").Value = -3389.8333
$ws.Range("H58").Value = 2591.375
$ws.Range("I58").Value = 1870.25
$ws.Range("J58").Value = 3312.5
$ws.Range("K58").Value = 1870.25
$ws.Range("L58").Value = 3312.5
$ws.Range("M58").Value = -1667.25
$ws.Range("N58").Value = -3718.5
$ws.Range("H132").Value = 2485.8333
$ws.Range("I132").Value = 2070.8948
$ws.Range("J132").Value = 4062.6
$ws.Range("K132").Value = 6212.6844
$ws.Range("L132").Value = 12187.8
$ws.Range("M132").Value = -3682.6844
$ws.Range("N132").Value = -17247.8
$ws.Range("H134").Value = 5600.963
$ws.Range("I134").Value = 6687.8096
$ws.Range("J134").Value = 1797
$ws.Range("K134").Value = 20063.4288
$ws.Range("L134").Value = 5391
$ws.Range("M134").Value = -17528.4288
$ws.Range("N134").Value = -10461
$ws.Range("H136").Value = 2591.375
$ws.Range("I136").Value = 1870.25
$ws.Range("J136").Value = 3312.5
$ws.Range("K136").Value = 5610.75
$ws.Range("L136").Value = 9937.5
$ws.Range("M136").Value = -3060.75
$ws.Range("N136").Value = -15037.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 468.32257
$ws.Range("I113").Value = 431.48486
$ws.Range("J113").Value = 510.2414
$ws.Range("K113").Value = 1294.45458
$ws.Range("L113").Value = 1530.7242
$ws.Range("M113").Value = 875.5454199999999
$ws.Range("N113").Value = -5870.724200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3359.8125
$ws.Range("J80").Value = 3645.7778
$ws.Range("L80").Value = 3645.7778
$ws.Range("N80").Value = -5641.7778
$ws.Range("H83").Value = 3359.8125
$ws.Range("J83").Value = 3645.7778
$ws.Range("L83").Value = 18228.889
$ws.Range("N83").Value = -28212.889
$ws.Range("H102").Value = 1696.2174
$ws.Range("I102").Value = 1545.1666
$ws.Range("K102").Value = 1545.1666
$ws.Range("M102").Value = 76.83339999999998
$ws.Range("H132").Value = 2269.9565
$ws.Range("I132").Value = 1911.7222
$ws.Range("J132").Value = 3559.6
$ws.Range("K132").Value = 5735.1666
$ws.Range("L132").Value = 10678.8
$ws.Range("M132").Value = -3205.1666
$ws.Range("N132").Value = -15738.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2030.5333
$ws.Range("I93").Value = 1711.4445
$ws.Range("J93").Value = 2509.1667
$ws.Range("K93").Value = 1711.4445
$ws.Range("L93").Value = 2509.1667
$ws.Range("M93").Value = -463.4445000000001
$ws.Range("N93").Value = -5005.1667
$ws.Range("H132").Value = 2471.0476
$ws.Range("I132").Value = 2053.7778
$ws.Range("J132").Value = 2784
$ws.Range("K132").Value = 6161.3334
$ws.Range("L132").Value = 8352
$ws.Range("M132").Value = -3631.3334
$ws.Range("N132").Value = -13412
$ws.Range("H136").Value = 1621.5834
$ws.Range("I136").Value = 1247.8235
$ws.Range("K136").Value = 3743.4705
$ws.Range("M136").Value = -1193.4705

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9976.916999999999
$ws.Range("I81").Value = 33973.668
$ws.Range("J81").Value = 1978
$ws.Range("K81").Value = 67947.336
$ws.Range("L81").Value = 3956
$ws.Range("M81").Value = -66886.336
$ws.Range("N81").Value = -6078
$ws.Range("H84").Value = 9976.916999999999
$ws.Range("I84").Value = 33973.668
$ws.Range("J84").Value = 1978
$ws.Range("K84").Value = 339736.68
$ws.Range("L84").Value = 19780
$ws.Range("M84").Value = -334432.68
$ws.Range("N84").Value = -30388
$ws.Range("H132").Value = 2827.5518
$ws.Range("I132").Value = 2762.4285
$ws.Range("J132").Value = 2998.5
$ws.Range("K132").Value = 8287.2855
$ws.Range("L132").Value = 8995.5
$ws.Range("M132").Value = -5757.2855
$ws.Range("N132").Value = -14055.5
$ws.Range("H136").Value = 1493.925
$ws.Range("I136").Value = 1262.2
$ws.Range("J136").Value = 3116
$ws.Range("K136").Value = 3786.6
$ws.Range("L136").Value = 9348
$ws.Range("M136").Value = -1236.6
$ws.Range("N136").Value = -14448
